$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:D18")
$rng.Sort($ws.Range("A2:A18"), 1)
